$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "group1"
$ws.Range("B1").Value = "group2"
$ws.Range("C1").Value = "meandiff"
$ws.Range("D1").Value = "p-adj"
$ws.Range("E1").Value = "lower"
$ws.Range("F1").Value = "upper"
$ws.Range("G1").Value = "reject"

$ws.Range("A1").Copy() | Out-Null
$ws.Range("B1:G1").PasteSpecial(-4122) | Out-Null

$ws.Range("A2").Value = "CSS"
$ws.Range("B2").Value = "Grassland"
$ws.Range("C2").Value = 0.452
$ws.Range("D2").Value = 0.001
$ws.Range("E2").Value = 0.271
$ws.Range("F2").Value = 0.6329
$ws.Range("G2").Value = $true
